$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Select()
$ws.Rows.Item(1).Insert()

$headerRange = $ws.Range("B2:G2")
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$bRange = $ws.Range("B3:B29")
$bRange.NumberFormat = "0.000"
$bRange.Borders.LineStyle = 1
$bRange.Borders.Weight = 2

$cgRange = $ws.Range("C3:G29")
$cgRange.Borders.LineStyle = 1
$cgRange.Borders.Weight = 2
